# ----------------------------------------------------------------------
# Reproduces the "adding remaining commits from last hw" edit:
#   - copies the two groups (A2:A49 / A50:A99) side-by-side into G/H so a
#     two-sample T.TEST can be redone against them in column I
#   - adds the T.TEST / SUM / AVERAGE helper formulas in D/E
#   - applies the black-font style that was used for the pasted H column
#   - restores the selection, column widths and page setup that came
#     along with the paste
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper formulas (top two groups split by Speaker, rows 2:49 / 50:99) ---
$ws.Range("D3").Formula = "=TTEST(A2:A49,A50:A99,1,3)"
$ws.Range("D6").Formula = "=SUM(A2:A49)"
$ws.Range("E6").Formula = "=AVERAGE(B2:B49)"
$ws.Range("D7").Formula = "=SUM(A50:A99)"
$ws.Range("E7").Formula = "=AVERAGE(A50:A99)"

# --- paste the two groups side by side into G (speaker 1) / H (speaker 2) ---
$ws.Range("A2:A49").Copy($ws.Range("G3:G50"))
$ws.Range("A50:A99").Copy($ws.Range("H3:H52"))

# the pasted-in "speaker 2" column got a plain black font applied
$ws.Range("H3:H52").Font.Color = 0

# redo the T.TEST on the freshly pasted columns
$ws.Range("I3").Formula = "=TTEST(G3:G50,H3:H52,1,3)"

# --- cosmetic bits that came along with the paste ---
$ws.Columns.Item(4).ColumnWidth = 11.33
$ws.Columns.Item(9).ColumnWidth = 11.33

$ws.Range("H3:H52").Select()

$ws.PageSetup.Orientation = 1
